# Apply the trading-results update described by the commit:
#   "Trade #60 closed at 2026-02-18 00:23:11 - unknown UNKNOWN +0.000%"
#
# Net effect of the diff:
#  - Summary sheet totals roll forward (one more trade closed)
#  - Strategy Status row for `momentum` rolls forward (+1 trade, new stats)
#  - Trade #88 (the open `momentum` trade) is closed out with an early exit,
#    reflected both on "All Trades" (row 89) and the per-strategy "momentum"
#    sheet (row 19)
#  - A brand-new open trade (#117, MarketMaking) is appended both to
#    "All Trades" (row 118) and the per-strategy "MarketMaking" sheet (row 38)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.14   # Current Capital
$summary.Range("B4").Value = 0.25      # Total P&L $
$summary.Range("B5").Value = 0.06      # Total P&L %
$summary.Range("B6").Value = 88        # Total Trades
$summary.Range("B7").Value = 44        # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.20999999999999   # Capital
$status.Range("D11").Value = 18                  # Trades
$status.Range("E11").Value = -0.79               # P&L $
$status.Range("F11").Value = -0.79               # P&L %
$status.Range("G11").Value = 22.22               # Win Rate %

# ---------------------------------------------------------------------------
# 3. All Trades sheet - close out Trade #88 (row 89) + append Trade #117 (row 118)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G89").Value = 0.97
$allTrades.Range("H89").Value = "CLOSED"
$allTrades.Range("I89").Value = 2.1053
$allTrades.Range("J89").Value = 0.02
$allTrades.Range("K89").Value = 99.20999999999999
$allTrades.Range("L89").Value = "early_exit"
$allTrades.Range("M89").Value = 0.12

$allTrades.Range("A118").Value = 117
$allTrades.Range("B118").Value = "'2026-02-18"
$allTrades.Range("C118").Value = "00:23:06"
$allTrades.Range("D118").Value = "MarketMaking"
$allTrades.Range("E118").Value = "DOWN"
$allTrades.Range("F118").Value = 0.95
$allTrades.Range("H118").Value = "OPEN"
$allTrades.Range("I118").Value = 0
$allTrades.Range("J118").Value = 0
$allTrades.Range("K118").Value = 99.410254715139
$allTrades.Range("M118").Value = 0
$allTrades.Range("N118").Value = 0
$allTrades.Range("O118").Value = 0
$allTrades.Range("P118").Value = 0.6
$allTrades.Range("Q118").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# 4. momentum sheet - mirror Trade #88 close-out (row 19)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("G19").Value = 0.97
$momentum.Range("H19").Value = "CLOSED"
$momentum.Range("I19").Value = 2.1053
$momentum.Range("J19").Value = 0.02
$momentum.Range("K19").Value = 99.20999999999999
$momentum.Range("P19").Value = "early_exit"
$momentum.Range("Q19").Value = 0.12

# ---------------------------------------------------------------------------
# 5. MarketMaking sheet - mirror the new Trade #117 (row 38)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Range("A38").Value = 117
$marketMaking.Range("B38").Value = "'2026-02-18"
$marketMaking.Range("C38").Value = "00:23:06"
$marketMaking.Range("D38").Value = "MarketMaking"
$marketMaking.Range("E38").Value = "DOWN"
$marketMaking.Range("F38").Value = 0.95
$marketMaking.Range("H38").Value = "OPEN"
$marketMaking.Range("I38").Value = 0
$marketMaking.Range("J38").Value = 0
$marketMaking.Range("K38").Value = 99.410254715139
$marketMaking.Range("L38").Value = 0
$marketMaking.Range("M38").Value = 0
$marketMaking.Range("N38").Value = 0.6
$marketMaking.Range("O38").Value = "Normal spread capture: 198 bps"
$marketMaking.Range("Q38").Value = 0
